$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-6) before rewriting full dataset (rows 2-8)
$ws.Range("A2:AQ6").ClearContents()

# Row 2
$ws.Range("A2").Value = 'Belgium'
$ws.Range("B2").Value = "'6"
$ws.Range("C2").Value = 'Healthcare Products'
$ws.Range("D2").Value = 0.0165
$ws.Range("E2").Value = -0.51
$ws.Range("G2").Value = 0.02679517762718002
$ws.Range("H2").Value = -0.2221485779022269
$ws.Range("I2").Value = -0.2755332262380862
$ws.Range("J2").Value = -0.252572124051579
$ws.Range("K2").Value = -101.838
$ws.Range("L2").Value = -0.2973485124982846
$ws.Range("M2").Value = 2.5542
$ws.Range("N2").Value = 0.001767135513598406
$ws.Range("O2").Value = -0.02508101101749838
$ws.Range("P2").Value = 2.5542
$ws.Range("Q2").Value = 0.001767135513598406
$ws.Range("R2").Value = -0.02508101101749838
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 320.028
$ws.Range("V2").Value = 0.2214129058593182
$ws.Range("W2").Value = -1.260598344352237
$ws.Range("X2").Value = 0.0595568855910163
$ws.Range("Y2").Value = -1.320155229943253
$ws.Range("Z2").Value = 1.18016459692212
$ws.Range("AA2").Value = -0.3543155040876804
$ws.Range("AB2").Value = 0.05850564141532216
$ws.Range("AC2").Value = -0.4126270324093066
$ws.Range("AD2").Value = 304.83
$ws.Range("AE2").Value = 0.2477402730171646
$ws.Range("AF2").Value = 305.0777402730172
$ws.Range("AG2").Value = -14.95025972698284
$ws.Range("AH2").Value = 0.1742835547631599
$ws.Range("AI2").Value = 0.5873117093794094
$ws.Range("AJ2").Value = -0.01045151312989206
$ws.Range("AK2").Value = -0.07496880552805386
$ws.Range("AL2").Value = 13.746
$ws.Range("AM2").Value = 11.706
$ws.Range("AN2").Value = -3.651314607414506
$ws.Range("AO2").Value = -6.877055143314418
$ws.Range("AP2").Value = 0.1790771962266615
$ws.Range("AQ2").Value = -8.075516828976593

# Row 3
$ws.Range("A3").Value = 'Belgium'
$ws.Range("B3").Value = 'Sequana Medical NV (ENXTBR:SEQUA)'
$ws.Range("C3").Value = 'Healthcare Products'
$ws.Range("G3").Value = -10.48461538461538
$ws.Range("H3").Value = -11.53846153846154
$ws.Range("I3").Value = -13.38461538461538
$ws.Range("J3").Value = -13.38461538461538
$ws.Range("K3").Value = -19.4
$ws.Range("L3").Value = -14.92307692307692
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 16.7
$ws.Range("V3").Value = 0.07330992098331869
$ws.Range("W3").Value = -1.993833504624871
$ws.Range("X3").Value = 0.0591445549118625
$ws.Range("Y3").Value = -2.052978059536734
$ws.Range("Z3").Value = -0.7602339181286554
$ws.Range("AA3").Value = 10.17543859649123
$ws.Range("AB3").Value = 0.05845471468679571
$ws.Range("AC3").Value = 10.11698388180444
$ws.Range("AD3").Value = 3.83
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 3.83
$ws.Range("AG3").Value = -12.87
$ws.Range("AH3").Value = 0.01653499114967837
$ws.Range("AI3").Value = 0.2635925671025465
$ws.Range("AJ3").Value = -0.05987996091750802
$ws.Range("AK3").Value = 5.930875576036867
$ws.Range("AL3").Value = 0.979
$ws.Range("AM3").Value = 0.941
$ws.Range("AN3").Value = -0.2201149425287356
$ws.Range("AO3").Value = -17.7732379979571
$ws.Range("AP3").Value = 0.7396551724137931
$ws.Range("AQ3").Value = -18.49096705632306

# Row 4
$ws.Range("A4").Value = 'Belgium'
$ws.Range("B4").Value = 'Ion Beam Applications SA (ENXTBR:IBAB)'
$ws.Range("C4").Value = 'Healthcare Products'
$ws.Range("D4").Value = 0.0165
$ws.Range("E4").Value = -0.51
$ws.Range("G4").Value = 0.1376051126807938
$ws.Range("H4").Value = 0.01281533804238143
$ws.Range("I4").Value = -0.0227716111671712
$ws.Range("J4").Value = -0.0113858055835856
$ws.Range("K4").Value = 0.995
$ws.Range("L4").Value = 0.003346787756474941
$ws.Range("M4").Value = 2.5542
$ws.Range("N4").Value = 0.00577089923181202
$ws.Range("O4").Value = 2.567035175879397
$ws.Range("P4").Value = 2.5542
$ws.Range("Q4").Value = 0.00577089923181202
$ws.Range("R4").Value = 2.567035175879397
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 107.9
$ws.Range("V4").Value = 0.2437867148666968
$ws.Range("W4").Value = 0.008868092691622103
$ws.Range("X4").Value = 0.06675045278625064
$ws.Range("Y4").Value = -0.05788236009462853
$ws.Range("Z4").Value = 1.529320987654321
$ws.Range("AA4").Value = -0.01741255144032922
$ws.Range("AB4").Value = 0.05738728587711735
$ws.Range("AC4").Value = -0.07479983731744658
$ws.Range("AD4").Value = 103.7
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 103.7
$ws.Range("AG4").Value = -4.200000000000003
$ws.Range("AH4").Value = 0.1898224418817499
$ws.Range("AI4").Value = 0.4841269841269842
$ws.Range("AJ4").Value = -0.009580291970802925
$ws.Range("AK4").Value = -0.03951081843838197
$ws.Range("AL4").Value = 2.16
$ws.Range("AM4").Value = 0.2300000000000002
$ws.Range("AN4").Value = -36.00694444444445
$ws.Range("AO4").Value = -3.134259259259259
$ws.Range("AP4").Value = 1.458333333333334
$ws.Range("AQ4").Value = -29.43478260869562

# Row 5
$ws.Range("A5").Value = 'Belgium'
$ws.Range("B5").Value = 'Remedent, Inc. (OTCPK:REMI)'
$ws.Range("C5").Value = 'Healthcare Products'
$ws.Range("D5").Value = -0.207
$ws.Range("G5").Value = -0.3885714285714285
$ws.Range("H5").Value = -0.3885714285714285
$ws.Range("I5").Value = -0.1671791830182426
$ws.Range("J5").Value = -0.1671791830182426
$ws.Range("K5").Value = -0.502
$ws.Range("L5").Value = -0.4780952380952381
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 0.128
$ws.Range("V5").Value = 0.04429065743944637
$ws.Range("W5").Value = -0.09940594059405941
$ws.Range("X5").Value = 0.05961772173097446
$ws.Range("Y5").Value = -0.1590236623250339
$ws.Range("Z5").Value = 0.2078100679596594
$ws.Range("AA5").Value = -0.03474151738446134
$ws.Range("AB5").Value = 0.05860056312930425
$ws.Range("AC5").Value = -0.0933420805137656
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0.08769071084577393
$ws.Range("AF5").Value = 0.08769071084577393
$ws.Range("AG5").Value = -0.04030928915422607
$ws.Range("AH5").Value = 0.02944923410828874
$ws.Range("AI5").Value = 0.01882707897318723
$ws.Range("AJ5").Value = -0.01414514529622847
$ws.Range("AK5").Value = -0.00889890540599396
$ws.Range("AL5").Value = 0.013
$ws.Range("AM5").Value = 0.013
$ws.Range("AN5").Value = 0
$ws.Range("AO5").Value = -18.76923076923077
$ws.Range("AP5").Value = 0.3698099922406061
$ws.Range("AQ5").Value = -18.76923076923077

# Row 6
$ws.Range("A6").Value = 'Belgium'
$ws.Range("B6").Value = 'Biocartis Group NV (ENXTBR:BCART)'
$ws.Range("C6").Value = 'Healthcare Products'
$ws.Range("D6").Value = 0.203
$ws.Range("G6").Value = -0.2827102803738318
$ws.Range("H6").Value = -1.336448598130841
$ws.Range("I6").Value = -1.453271028037383
$ws.Range("J6").Value = -1.453271028037383
$ws.Range("K6").Value = -74.2
$ws.Range("L6").Value = -1.733644859813084
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("U6").Value = 168.4
$ws.Range("V6").Value = 0.5145126794989306
$ws.Range("W6").Value = -0.527363184079602
$ws.Range("X6").Value = 0.07843075815019337
$ws.Range("Y6").Value = -0.6057939422297954
$ws.Range("Z6").Value = 0.4637053087757314
$ws.Range("AA6").Value = -0.6738894907908994
$ws.Range("AB6").Value = 0.05802249351394818
$ws.Range("AC6").Value = -0.7319119843048476
$ws.Range("AD6").Value = 186
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 186
$ws.Range("AG6").Value = 17.59999999999999
$ws.Range("AH6").Value = 0.3623611922852134
$ws.Range("AI6").Value = 0.7560975609756098
$ws.Range("AJ6").Value = 0.0510292838503914
$ws.Range("AK6").Value = 0.2268041237113401
$ws.Range("AL6").Value = 10.5
$ws.Range("AM6").Value = 10.5
$ws.Range("AN6").Value = -3.351351351351351
$ws.Range("AO6").Value = -5.923809523809524
$ws.Range("AP6").Value = -0.317117117117117
$ws.Range("AQ6").Value = -5.923809523809524

# Row 7
$ws.Range("A7").Value = 'Belgium'
$ws.Range("B7").Value = 'Nyxoah S.A. (ENXTBR:NYXH)'
$ws.Range("C7").Value = 'Healthcare Products'
$ws.Range("K7").Value = -8.61
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("U7").Value = 26.9
$ws.Range("V7").Value = 0.06305672761368963
$ws.Range("X7").Value = 0.05949604945105814
$ws.Range("Z7").Value = 0
$ws.Range("AA7").Value = -48.12890337172826
$ws.Range("AB7").Value = 0.05859565061552971
$ws.Range("AC7").Value = -48.18749902234379
$ws.Range("AD7").Value = 11.3
$ws.Range("AE7").Value = 0.1600495621713907
$ws.Range("AF7").Value = 11.46004956217139
$ws.Range("AG7").Value = -15.43995043782861
$ws.Range("AH7").Value = 0.02616091007072064
$ws.Range("AI7").Value = 0.2860717769304282
$ws.Range("AJ7").Value = -0.03755216600997596
$ws.Range("AK7").Value = -1.173244094931891
$ws.Range("AL7").Value = 0.091
$ws.Range("AM7").Value = 0.019
$ws.Range("AN7").Value = -1.504460125149781
$ws.Range("AO7").Value = -85.71428571428571
$ws.Range("AP7").Value = 2.055645112212569
$ws.Range("AQ7").Value = -410.5263157894736

# Row 8
$ws.Range("A8").Value = 'Belgium'
$ws.Range("B8").Value = 'Metrics in Balance N.V. (ENXTPA:MLMIB)'
$ws.Range("C8").Value = 'Healthcare Products'
$ws.Range("G8").Value = -2.297297297297297
$ws.Range("H8").Value = -2.297297297297297
$ws.Range("I8").Value = -3.189189189189189
$ws.Range("J8").Value = -3.189189189189189
$ws.Range("K8").Value = -0.121
$ws.Range("L8").Value = -3.27027027027027
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("U8").Value = 0
$ws.Range("V8").Value = 0
$ws.Range("X8").Value = 0.0585565681438486
$ws.Range("AB8").Value = 0.0585565681438486
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AJ8").Value = 0
$ws.Range("AL8").Value = 0.003
$ws.Range("AM8").Value = 0.003
$ws.Range("AN8").Value = 0
$ws.Range("AO8").Value = -39.33333333333333
$ws.Range("AP8").Value = 0
$ws.Range("AQ8").Value = -39.33333333333333
